# sdl-7869: re-add test cases to the getDataGraphQL sheet (rows 4-9),
# and restore the "getDataGraphQL" tab as the active/selected sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("getDataGraphQL")

# --- Copy the existing bordered cell format down into rows 4-9, matching
# --- the column A/B "bordered" style used throughout the sheet.
$ws1.Range("A3").Copy()
$ws1.Range("A4:F9").PasteSpecial(-4122)

# Columns D:F (and C for rows 5-9) in the new rows were filled without the
# border formatting used by columns A/B - drop the border there so the new
# rows reuse the sheet's existing "no border" style instead of creating a
# brand new one.
$ws1.Range("D4:F9").Borders.LineStyle = 0
$ws1.Range("C5:C9").Borders.LineStyle = 0

# --- Row 4: databrain-getDataGraphQL-test-3
$ws1.Range("A4").Value = "databrain-getDataGraphQL-test-3"
$ws1.Range("B4").Value = "good request, data retrieved"
$ws1.Range("C4").Value = "{Equipment {type sourceId description _isPartOf id name}}"
$ws1.Range("D4").Value = 200
$ws1.Range("E4").Value = 100000
$ws1.Range("F4").Value = "Successfully"

# --- Row 5: databrain-getDataGraphQL-test-4
$ws1.Range("A5").Value = "databrain-getDataGraphQL-test-4"
$ws1.Range("B5").Value = "good request, data retrieved"
$ws1.Range("C5").Value = "{Location {type sourceId description _isPartOf id name}}"
$ws1.Range("D5").Value = 200
$ws1.Range("E5").Value = 100000
$ws1.Range("F5").Value = "Successfully"

# --- Row 6: databrain-getDataGraphQL-test-5
$ws1.Range("A6").Value = "databrain-getDataGraphQL-test-5"
$ws1.Range("B6").Value = "good request, data retrieved"
$ws1.Range("C6").Value = "{Point {sourceId sourceSystem type id description _isPointOf name}}"
$ws1.Range("D6").Value = 200
$ws1.Range("E6").Value = 100000
$ws1.Range("F6").Value = "Successfully"

# --- Row 7: databrain-getDataGraphQL-test-6
$ws1.Range("A7").Value = "databrain-getDataGraphQL-test-6"
$ws1.Range("B7").Value = "good request, data retrieved"
$ws1.Range("C7").Value = "{Sensor {sourceId sourceSystem type id description _isPointOf name}}"
$ws1.Range("D7").Value = 200
$ws1.Range("E7").Value = 100000
$ws1.Range("F7").Value = "Successfully"

# --- Row 8: databrain-getDataGraphQL-test-11
$ws1.Range("A8").Value = "databrain-getDataGraphQL-test-11"
$ws1.Range("B8").Value = "good request, data retrieved"
$ws1.Range("C8").Value = "{Humidity_Sensor {sourceId sourceSystem type id description _isPointOf name isPointOf_Thermostat {type sourceId description _isPartOf id name}}}"
$ws1.Range("D8").Value = 200
$ws1.Range("E8").Value = 100000
$ws1.Range("F8").Value = "Successfully"

# --- Row 9: databrain-getDataGraphQL-test-12
$ws1.Range("A9").Value = "databrain-getDataGraphQL-test-12"
$ws1.Range("B9").Value = "good request, data retrieved"
$ws1.Range("C9").Value = "{Room {type sourceId description _isPartOf id name isPartOf_Floor {type sourceId description _isPartOf id name}}}"
$ws1.Range("D9").Value = 200
$ws1.Range("E9").Value = 100000
$ws1.Range("F9").Value = "Successfully"

# --- Row heights: the new/refreshed rows render at 26.4pt.
$ws1.Rows.Item(1).RowHeight = 26.4
$ws1.Rows.Item(2).RowHeight = 26.4
$ws1.Rows.Item(3).RowHeight = 26.4
$ws1.Range("A4:F9").RowHeight = 26.4

$ws2 = $wb.Worksheets.Item("queryUpdateFromApiEngine")
$ws2.Rows.Item(1).RowHeight = 26.4
$ws2.Rows.Item(2).RowHeight = 79.2
$ws2.Rows.Item(3).RowHeight = 79.2
$ws2.Rows.Item(4).RowHeight = 52.8
$ws2.Rows.Item(5).RowHeight = 224.4

$ws3 = $wb.Worksheets.Item("queryPostgresqlData")
$ws3.Rows.Item(1).RowHeight = 26.4

# --- Re-select getDataGraphQL as the active sheet/cell (it had been
# --- queryPostgresqlData); move the cursor to C17 to match.
$ws1.Activate()
$ws1.Range("C17").Select()

Write-Host "done"
